# -------------------------------------------------------------------------
# Avenant_rdm_026.docx : the "Référence de l'étude" placeholder changes
# from the Jinja-style token {{ref_etude}} to {{etude.ref}}.
#
# In the canonical OOXML this shows up as the single run that holds the
# text "ref_etude" being split into two runs - "etude" and ".ref" - each
# wrapped by a <w:proofErr w:type="spellStart"/> / <w:proofErr
# w:type="spellEnd"/> pair (the artefact Word leaves behind when the
# in-place edit is re-checked by the spell checker), while every other
# run in the paragraph (the two "{" runs, the two "}" runs, the label,
# etc.) stays untouched.
# -------------------------------------------------------------------------

$d = $word.ActiveDocument

# Find the placeholder so the edit is anchored to its actual location
# rather than a hard-coded character offset.
$rng = $d.Content
$find = $rng.Find
$find.Text = "ref_etude"
$find.Execute() | Out-Null

if (-not $find.Found) {
    throw "Could not find the 'ref_etude' placeholder in the document."
}

$hit = $d.Range($rng.Start, $rng.End)
$para = $hit.Paragraphs(1)
$pRange = $para.Range
$pStart = $pRange.Start
$pEnd = $pRange.End

# The paragraph's character style id ("Rfrencelgre") - reapplied below,
# because InsertXML (unlike normal typing) does not resolve rStyle
# references on its own.
$styleId = $para.Range.Style

# Full replacement for the paragraph: identical to the original except
# that the run containing "ref_etude" is replaced by two runs - "etude"
# and ".ref" - bracketed by proofErr spell-check markers.
$newParagraphXml = @'
<w:p w14:paraId="7AAA51ED" w14:textId="5FDABD0F" w:rsidR="002C69D3" w:rsidRPr="00347DE4" w:rsidRDefault="002C69D3" w:rsidP="002C69D3"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="right"/><w:rPr><w:rStyle w:val="Rfrencelgre"/><w:rFonts w:ascii="Quicksand" w:hAnsi="Quicksand"/><w:color w:val="7D92DF"/></w:rPr></w:pPr><w:r w:rsidRPr="00347DE4"><w:rPr><w:rStyle w:val="Rfrencelgre"/><w:rFonts w:ascii="Quicksand" w:hAnsi="Quicksand"/><w:color w:val="7D92DF"/></w:rPr><w:t>Référence de l’étude :</w:t></w:r><w:r w:rsidR="00967F4A"><w:rPr><w:rStyle w:val="Rfrencelgre"/><w:rFonts w:ascii="Quicksand" w:hAnsi="Quicksand"/><w:color w:val="7D92DF"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00347DE4"><w:rPr><w:rStyle w:val="Rfrencelgre"/><w:rFonts w:ascii="Quicksand" w:hAnsi="Quicksand"/><w:color w:val="7D92DF"/></w:rPr><w:t>{</w:t></w:r><w:r w:rsidR="004F0A80"><w:rPr><w:rStyle w:val="Rfrencelgre"/><w:rFonts w:ascii="Quicksand" w:hAnsi="Quicksand"/><w:color w:val="7D92DF"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="005E5F5E"><w:rPr><w:rStyle w:val="Rfrencelgre"/><w:rFonts w:ascii="Quicksand" w:hAnsi="Quicksand"/><w:color w:val="7D92DF"/></w:rPr><w:t>etude</w:t></w:r><w:r w:rsidR="005E5F5E"><w:rPr><w:rStyle w:val="Rfrencelgre"/><w:rFonts w:ascii="Quicksand" w:hAnsi="Quicksand"/><w:color w:val="7D92DF"/></w:rPr><w:t>.ref</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004F0A80"><w:rPr><w:rStyle w:val="Rfrencelgre"/><w:rFonts w:ascii="Quicksand" w:hAnsi="Quicksand"/><w:color w:val="7D92DF"/></w:rPr><w:t>}</w:t></w:r><w:r w:rsidRPr="00347DE4"><w:rPr><w:rStyle w:val="Rfrencelgre"/><w:rFonts w:ascii="Quicksand" w:hAnsi="Quicksand"/><w:color w:val="7D92DF"/></w:rPr><w:t>}</w:t></w:r></w:p>
'@

$pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part" pkg:contentType="application/xml"><pkg:xmlData xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $newParagraphXml + '</pkg:xmlData></pkg:part></pkg:package>'

$pRange.InsertXML($pkg) | Out-Null

# InsertXML does not preserve the rStyle reference on the runs it
# creates, so reapply the paragraph's character style across the whole
# (just replaced) paragraph to restore it, matching the original
# formatting exactly.
$fixup = $d.Range($pStart, $pEnd)
$fixup.Style = $styleId
